# Import product and variant
# Adds a "Taxons" column and a "Promotionable" column to the Main sheet,
# updates the existing product row's Detail/Cost/Discontinue-date values,
# adds a second product row (Winter 2021 Tshirt), and moves the two header
# comments to follow the shifted "Sale Price" / "Suggested Retail Price"
# columns. The Variant sheet's formulas/shared strings follow automatically.

$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("Main")
$variant = $wb.Worksheets.Item("Variant")

# --- 1. Insert the two new columns -----------------------------------
# "Taxons" goes in before the old column E ("Option Types").
$main.Range("E1").EntireColumn.Insert()
# "Promotionable" goes in before the (now shifted) "Available Date" column.
$main.Range("N1").EntireColumn.Insert()

# --- 2. New column headers --------------------------------------------
$main.Range("E1").Value = "Taxons"
$main.Range("N1").Value = "Promotionable"

# --- 3. Move the header comments ---------------------------------------
# They used to sit on J1/K1 (Sale Price / Suggested Retail Price); those
# columns are now K1/L1 after the "Taxons" insert.
$main.Range("J1").Comment.Delete()
$main.Range("K1").Comment.Delete()
$main.Range("K1").AddComment("Master Price")
$main.Range("L1").AddComment("Compared at Price")

# --- 4. Update the existing product row (row 2) -------------------------
$main.Range("B2").Value = "Shirt"
$main.Range("J2").Value = "Check out this tshirt as part of our summer 2021 lookbook."
$main.Range("M2").Value = "USD 10"
$main.Range("P2").Value = "2022-02-30"

# --- 5. Add the new product row (row 3) ---------------------------------
$main.Range("A3").Value = "Zando Fashion"
$main.Range("B3").Value = "Shirt"
$main.Range("C3").Value = "Shipping by VTENH"
$main.Range("D3").Value = "VAT Incl."
$main.Range("E3").Value = "Womenware"
$main.Range("F3").Value = "Color"
$main.Range("G3").Value = "Winter 2021 Tshirt"
$main.Range("H3").Value = "WTST202121"
$main.Range("I3").Value = "Trendy Tshirt for your summer need. Check out this tshirt as part of our summer 2021 lookbook."
$main.Range("J3").Value = "Awesome out this tshirt as part of our summer 2021 lookbook."
$main.Range("K3").Value = "USD 15"
$main.Range("L3").Value = "USD 18"
$main.Range("M3").Value = "USD 12"
$main.Range("N3").Value = "yes"
$main.Range("O3").Value = "2021-02-30"
$main.Range("Q3").Value = "Winter 2021 Tshirt"
$main.Range("R3").Value = "Beautiful Winter 2021 Tshirt"
$main.Range("T3").Value = "Cotton"
$main.Range("U3").Value = "Nike"
$main.Range("V3").Value = "ZFMST202121-1.jpg"
$main.Range("W3").Value = "ZFMST202121-2.jpg"
$main.Range("X3").Value = "ZFMST202121-3.jpg"
$main.Range("Y3").Value = "ZFMST202121-4.jpg"
